# Applies the "Updated imputation procedure, first steps of imputed case analysis"
# changes to the Comp_comb_plant.xlsx data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Column B ("Method"): every "Repellent (...)" label (rows 14-148) becomes
#    the single unified label "Repellents".
# ---------------------------------------------------------------------------
for ($r = 14; $r -le 148; $r++) {
    $ws.Cells.Item($r, 2).Value = "Repellents"
}

# ---------------------------------------------------------------------------
# 2) Column D ("Animal"): specific animal labels get collapsed into broader
#    categories ("Other ungulates" / "Other").
# ---------------------------------------------------------------------------
$dRanges = @(
    @{Start=2;   End=13;  Value="Other ungulates"},
    @{Start=21;  End=24;  Value="Other"},
    @{Start=31;  End=45;  Value="Other"},
    @{Start=52;  End=68;  Value="Other"},
    @{Start=69;  End=69;  Value="Other ungulates"},
    @{Start=74;  End=78;  Value="Other ungulates"},
    @{Start=79;  End=130; Value="Other"},
    @{Start=131; End=148; Value="Other ungulates"}
)

foreach ($range in $dRanges) {
    for ($r = $range.Start; $r -le $range.End; $r++) {
        $ws.Cells.Item($r, 4).Value = $range.Value
    }
}

# ---------------------------------------------------------------------------
# 3) Columns O, P, R, S (Con_mean, Con_SD, Trt_mean, Trt_SD) for rows 46-78
#    were recomputed under the new (logit-based) imputation procedure.
# ---------------------------------------------------------------------------
$oprsValues = @(
    @{Row=46; O=0.6645915469470783; P=2.456123440944749; R=-1.074913191092807; S=5.011375609674042},
    @{Row=47; O=0.6645915469470783; P=2.456123440944749; R=-1.564702018592708; S=5.438386561619762},
    @{Row=48; O=0.6645915469470783; P=2.456123440944749; R=-2.256357432280105; S=5.616825552931346},
    @{Row=49; O=0.9618533515527208; P=2.346397226232656; R=-1.431074030942754; S=5.755318862234189},
    @{Row=50; O=0.9618533515527208; P=2.346397226232656; R=-2.059482652175352; S=6.725414908905554},
    @{Row=51; O=0.9618533515527208; P=2.346397226232656; R=0.352031754008015; S=2.604696105940015},
    @{Row=52; O=0.2896332925830415; P=10.11114316401136; R=-1.340925756090039; S=15.77425434967047},
    @{Row=53; O=0.3544336280957182; P=12.18463279010214; R=-1.386294361119889; S=20.88678646331025},
    @{Row=54; O=0.3880300462860632; P=8.330703146777585; R=-1.363988603605591; S=14.25121975864779},
    @{Row=55; O=0.8506535682341779; P=7.981138506182281; R=-2.959364629383103; S=30.08969389127185},
    @{Row=56; O=1.324221932477512; P=4.549834668506687; R=-3.733693469903739; S=22.93131799524272},
    @{Row=57; O=-0.8839553512312189; P=9.90930057676136; R=-0.2336148511815058; S=7.112155534885336},
    @{Row=58; O=0.1945916493857688; P=4.424604026796255; R=-1.058871960018595; S=6.457238021401541},
    @{Row=59; O=0.1801261662305197; P=3.477192159559847; R=-1.222226244635286; S=5.842071671790914},
    @{Row=60; O=0.141775462399605; P=3.708141601196324; R=-1.203972804325934; S=5.177284281908125},
    @{Row=61; O=1.145865173518654; P=2.626016974719664; R=-2.932218867831669; S=10.16864845566765},
    @{Row=62; O=4.248495242049376; P=1.737240847946985; R=0.4734380916943716; S=1.977676762508595},
    @{Row=63; O=1.475906519809576; P=1.894172353707069; R=3.600731067337229; S=1.319670139537233},
    @{Row=64; O=1.885691289062258; P=2.597245857243066; R=0.1078889620111851; S=3.434933587141446},
    @{Row=65; O=1.880706955244552; P=2.508225041019831; R=0.1047109540890318; S=3.428807907634629},
    @{Row=66; O=2.74798291143758; P=2.028112935206314; R=0.03344793406753969; S=2.642615857046959},
    @{Row=67; O=3.919991175077309; P=0.6736852996662777; R=-1.751485570090117; S=2.229765426227599},
    @{Row=68; O=2.094479419044374; P=0.8623590898619642; R=4.307437777682794; S=0.6553181618869695},
    @{Row=69; O=0.4843364116856495; P=3.570900437363647; R=-0.1271334181448639; S=5.357937202821788},
    @{Row=70; O=0.2258945846860995; P=5.246596843224606; R=-0.2905351414341878; S=5.433162787804286},
    @{Row=71; O=-0.6270700899881577; P=1.084061345602135; R=-1.258873445792383; S=1.04752469561917},
    @{Row=72; O=-0.5038244283070068; P=1.099153454338057; R=-1.579947483923394; S=1.04376967914519},
    @{Row=73; O=-0.560189346884326; P=1.08991857504101; R=-1.432556078935963; S=1.042468625030766},
    @{Row=74; O=-0.6378845018848964; P=1.083818549022401; R=-2.052949926232858; S=1.031342514579477},
    @{Row=75; O=0.01635724202477102; P=1.171163262864323; R=-3.363457053679302; S=1.029854453722139},
    @{Row=76; O=-0.707680695176113; P=1.07844923740431; R=-3.363457053679302; S=1.029854453722139},
    @{Row=77; O=-0.0262310121488693; P=1.162737379669453; R=-0.06067434461186064; S=1.155336569909553},
    @{Row=78; O=-0.0262310121488693; P=1.162737379669453; R=-2.161314076030447; S=1.025000810318442}
)

foreach ($entry in $oprsValues) {
    $r = $entry.Row
    $ws.Cells.Item($r, 15).Value = $entry.O
    $ws.Cells.Item($r, 16).Value = $entry.P
    $ws.Cells.Item($r, 18).Value = $entry.R
    $ws.Cells.Item($r, 19).Value = $entry.S
}
